$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data for rows 2-11 (columns A-G)
$data = @(
    @(29505, "Rodrigo Porto",   "Vendas",                 "Doenca",             4, 45093, 6091.09),
    @(3592,  "Enzo Cavalcanti", "Vendas",                 "Consulta medica",    4, 45103, 2238.73),
    @(72260, "Zoe Casa Grande", "Atendimento ao Cliente", "Consulta medica",    1, 45091, 4212.95),
    @(99335, "Caio Farias",     "Operacoes",               "Viagem de negocios", 7, 45090, 3659.44),
    @(41419, "Julia Dias",      "Financeiro",              "Viagem de negocios", 7, 45104, 2071.02),
    @(36391, "Sophia Rocha",    "Engenharia",               "Problemas pessoais", 8, 45084, 5778.54),
    @(65809, "Danilo Ribeiro",  "TI",                       "Problemas pessoais", 3, 45105, 8286.11),
    @(88283, "Juan Nogueira",   "Atendimento ao Cliente",   "Outros",             7, 45087, 2360.61),
    @(74365, "Thomas Pimenta",  "Juridico",                 "Viagem de negocios", 7, 45082, 8471.49),
    @(77209, "Dom Mendes",      "Recursos Humanos",         "Outros",             3, 45105, 7608.46)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
